$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The naive component forecaster bug fix shifts every existing data row
# (rows 2-11) down by one row, dropping the previous last row (old row 11,
# N=9/Q9) off the bottom, and inserts freshly computed error-metric values
# into the new top data row (row 2, N=19/Q0).
#
# Column A (the Q-labels) is untouched by the fix, so only columns B:G are
# rewritten for rows 2-11.

$newData = @(
  @(0.02072117565895826, 0.5679342762134251, 0.6583304034161481, 0.8113756241200176, 0.8333373229775733, 19),
  @(0.04612272632165332, 0.6122205438342831, 0.7643369223465524, 0.8742636457880153, 0.8983570984052854, 18),
  @(-0.05484937042860891, 0.5505006181463401, 0.631985967850102, 0.7949754510990273, 0.8174892103658266, 17),
  @(0.0454318984618024, 0.4533283699478702, 0.4371073281485932, 0.6611409291131454, 0.6812093330742512, 16),
  @(0.1619522260889591, 0.3783700883078635, 0.2480121053477303, 0.4980081378328373, 0.4874682207362279, 15),
  @(0.1471867289230871, 0.4060743056329635, 0.2139495689888607, 0.4625468289685496, 0.455056961899273, 14),
  @(0.169821430505846, 0.3504591757809973, 0.229570748894733, 0.4791354181176058, 0.4663250474243488, 13),
  @(0.1844780524147233, 0.3492283006690807, 0.2072969125451589, 0.4552987069443081, 0.4347598217479242, 12),
  @(0.2035825112292451, 0.3072747587892817, 0.2238621271029474, 0.4731407053963413, 0.4479485651801955, 11),
  @(0.1565263708236955, 0.2762321533226782, 0.1325312503702182, 0.3640484176180666, 0.3464594662951709, 10)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
  $row = 2 + $i
  $vals = $newData[$i]
  $ws.Cells.Item($row, 2).Value = $vals[0]
  $ws.Cells.Item($row, 3).Value = $vals[1]
  $ws.Cells.Item($row, 4).Value = $vals[2]
  $ws.Cells.Item($row, 5).Value = $vals[3]
  $ws.Cells.Item($row, 6).Value = $vals[4]
  $ws.Cells.Item($row, 7).Value = $vals[5]
}
